$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

$ws.Cells.Item(2, 2).Value = 12.2618424620122
$ws.Cells.Item(2, 3).Value = 4.650237032091947
$ws.Cells.Item(2, 4).Value = 15.09202394598852
$ws.Cells.Item(2, 5).Value = 16.52090033994226
$ws.Cells.Item(2, 7).Value = 3.699854853063714
$ws.Cells.Item(2, 10).Value = 9.443349102603365
$ws.Cells.Item(2, 11).Value = 11.7255274943759
$ws.Cells.Item(2, 15).Value = 31.05843703583859

$ws.Cells.Item(3, 2).Value = 12.05188601409519
$ws.Cells.Item(3, 3).Value = 4.490052121096118
$ws.Cells.Item(3, 4).Value = 15.03404296433571
$ws.Cells.Item(3, 5).Value = 16.46339042045647
$ws.Cells.Item(3, 7).Value = 3.702265505868187
$ws.Cells.Item(3, 10).Value = 9.452115687963095
$ws.Cells.Item(3, 11).Value = 11.59089284022762
$ws.Cells.Item(3, 15).Value = 31.1087812818225

$ws.Cells.Item(4, 2).Value = 11.92404742369658
$ws.Cells.Item(4, 3).Value = 4.389631598517343
$ws.Cells.Item(4, 4).Value = 15.00177278806841
$ws.Cells.Item(4, 5).Value = 16.43175411890232
$ws.Cells.Item(4, 7).Value = 3.703823486022263
$ws.Cells.Item(4, 10).Value = 9.458926873260094
$ws.Cells.Item(4, 11).Value = 11.50995957591305
$ws.Cells.Item(4, 15).Value = 31.14579619755863

$ws.Cells.Item(5, 2).Value = 11.87229510258514
$ws.Cells.Item(5, 3).Value = 4.348257272616038
$ws.Cells.Item(5, 4).Value = 14.98946943966159
$ws.Cells.Item(5, 5).Value = 16.41979522586885
$ws.Cells.Item(5, 7).Value = 3.704478011973797
$ws.Cells.Item(5, 10).Value = 9.462061795475929
$ws.Cells.Item(5, 11).Value = 11.47745172940404
$ws.Cells.Item(5, 15).Value = 31.16241149343617

$ws.Cells.Item(6, 2).Value = 11.86372450434459
$ws.Cells.Item(6, 3).Value = 4.341361874508777
$ws.Cells.Item(6, 4).Value = 14.98747789401924
$ws.Cells.Item(6, 5).Value = 16.41786606375245
$ws.Cells.Item(6, 7).Value = 3.704587883424918
$ws.Cells.Item(6, 10).Value = 9.462604050042982
$ws.Cells.Item(6, 11).Value = 11.47208347585145
$ws.Cells.Item(6, 15).Value = 31.16526284070734

$ws.Cells.Item(7, 2).Value = 11.92334798930636
$ws.Cells.Item(7, 3).Value = 4.389075347187974
$ws.Cells.Item(7, 4).Value = 15.00160341945065
$ws.Cells.Item(7, 5).Value = 16.43158904764297
$ws.Cells.Item(7, 7).Value = 3.703832233594198
$ws.Cells.Item(7, 10).Value = 9.458967697001231
$ws.Cells.Item(7, 11).Value = 11.50951919940024
$ws.Cells.Item(7, 15).Value = 31.14601408130514

$ws.Cells.Item(8, 2).Value = 12.18926963805131
$ws.Cells.Item(8, 3).Value = 4.595472794945478
$ws.Cells.Item(8, 4).Value = 15.07134693153928
$ws.Cells.Item(8, 5).Value = 16.50031392187938
$ws.Cells.Item(8, 7).Value = 3.700669932121581
$ws.Cells.Item(8, 10).Value = 9.446075428486617
$ws.Cells.Item(8, 11).Value = 11.67876844179013
$ws.Cells.Item(8, 15).Value = 31.07452708449413

$ws.Cells.Item(9, 2).Value = 12.71611889554351
$ws.Cells.Item(9, 3).Value = 4.98115709385711
$ws.Cells.Item(9, 4).Value = 15.23407189939706
$ws.Cells.Item(9, 5).Value = 16.6637966851863
$ws.Cells.Item(9, 7).Value = 3.695083221568555
$ws.Cells.Item(9, 10).Value = 9.432121240747952
$ws.Cells.Item(9, 11).Value = 12.02274077562158
$ws.Cells.Item(9, 15).Value = 30.98290499681257

$ws.Cells.Item(10, 2).Value = 13.10233241754726
$ws.Cells.Item(10, 3).Value = 5.249741911032868
$ws.Cells.Item(10, 4).Value = 15.36876007266436
$ws.Cells.Item(10, 5).Value = 16.80074597399715
$ws.Cells.Item(10, 7).Value = 3.691349127095673
$ws.Cells.Item(10, 10).Value = 9.428762698924785
$ws.Cells.Item(10, 11).Value = 12.280546478235
$ws.Cells.Item(10, 15).Value = 30.94536448881258

$ws.Cells.Item(11, 2).Value = 13.27700420332106
$ws.Cells.Item(11, 3).Value = 5.368164355798664
$ws.Cells.Item(11, 4).Value = 15.43315286358379
$ws.Cells.Item(11, 5).Value = 16.86654137685085
$ws.Cells.Item(11, 7).Value = 3.68972994160504
$ws.Cells.Item(11, 10).Value = 9.428727568816409
$ws.Cells.Item(11, 11).Value = 12.39845040331282
$ws.Cells.Item(11, 15).Value = 30.93477982722794

$ws.Cells.Item(12, 2).Value = 13.34293485438944
$ws.Cells.Item(12, 3).Value = 5.412428208145581
$ws.Cells.Item(12, 4).Value = 15.45797036476537
$ws.Cells.Item(12, 5).Value = 16.8919439904057
$ws.Cells.Item(12, 7).Value = 3.689128158492463
$ws.Cells.Item(12, 10).Value = 9.428928393429221
$ws.Cells.Item(12, 11).Value = 12.44314767121283
$ws.Cells.Item(12, 15).Value = 30.93170688485388

$ws.Cells.Item(13, 2).Value = 13.3287460548379
$ws.Cells.Item(13, 3).Value = 5.402921580282323
$ws.Cells.Item(13, 4).Value = 15.45260645592268
$ws.Cells.Item(13, 5).Value = 16.8864516644649
$ws.Cells.Item(13, 7).Value = 3.689257258652908
$ws.Cells.Item(13, 10).Value = 9.428875628083786
$ws.Cells.Item(13, 11).Value = 12.43351974628589
$ws.Cells.Item(13, 15).Value = 30.93232708222398

$ws.Cells.Item(14, 2).Value = 13.2824330162933
$ws.Cells.Item(14, 3).Value = 5.371817807422059
$ws.Cells.Item(14, 4).Value = 15.43518602177715
$ws.Cells.Item(14, 5).Value = 16.86862159006938
$ws.Cells.Item(14, 7).Value = 3.689680205075462
$ws.Cells.Item(14, 10).Value = 9.428739802630485
$ws.Cells.Item(14, 11).Value = 12.40212690986348
$ws.Cells.Item(14, 15).Value = 30.93450826324984

$ws.Cells.Item(15, 2).Value = 13.2540351599013
$ws.Cells.Item(15, 3).Value = 5.352689197245051
$ws.Cells.Item(15, 4).Value = 15.42457145991293
$ws.Cells.Item(15, 5).Value = 16.85776314499051
$ws.Cells.Item(15, 7).Value = 3.689940750530615
$ws.Cells.Item(15, 10).Value = 9.428684474435114
$ws.Cells.Item(15, 11).Value = 12.38290316065529
$ws.Cells.Item(15, 15).Value = 30.93596613521132

$ws.Cells.Item(16, 2).Value = 13.09089088725948
$ws.Cells.Item(16, 3).Value = 5.241923594583382
$ws.Cells.Item(16, 4).Value = 15.36461336354422
$ws.Cells.Item(16, 5).Value = 16.79651514515664
$ws.Cells.Item(16, 7).Value = 3.691456538664212
$ws.Cells.Item(16, 10).Value = 9.428794989109349
$ws.Cells.Item(16, 11).Value = 12.2728503136416
$ws.Cells.Item(16, 15).Value = 30.94618701763048

$ws.Cells.Item(17, 2).Value = 12.99049743819496
$ws.Cells.Item(17, 3).Value = 5.172980109278705
$ws.Cells.Item(17, 4).Value = 15.32861986008684
$ws.Cells.Item(17, 5).Value = 16.75982628305076
$ws.Cells.Item(17, 7).Value = 3.692406737085936
$ws.Cells.Item(17, 10).Value = 9.429244758238447
$ws.Cells.Item(17, 11).Value = 12.20546796877192
$ws.Cells.Item(17, 15).Value = 30.95412141017975

$ws.Cells.Item(18, 2).Value = 12.93266174175565
$ws.Cells.Item(18, 3).Value = 5.13297396113845
$ws.Cells.Item(18, 4).Value = 15.30821223889573
$ws.Cells.Item(18, 5).Value = 16.73905390790978
$ws.Cells.Item(18, 7).Value = 3.692960750050306
$ws.Cells.Item(18, 10).Value = 9.429643989007081
$ws.Cells.Item(18, 11).Value = 12.16677281647076
$ws.Cells.Item(18, 15).Value = 30.95929611882149

$ws.Cells.Item(19, 2).Value = 12.91306585975029
$ws.Cells.Item(19, 3).Value = 5.11936945623779
$ws.Cells.Item(19, 4).Value = 15.3013536809633
$ws.Cells.Item(19, 5).Value = 16.73207788876006
$ws.Cells.Item(19, 7).Value = 3.693149616559283
$ws.Cells.Item(19, 10).Value = 9.429803313148479
$ws.Cells.Item(19, 11).Value = 12.15368309179976
$ws.Cells.Item(19, 15).Value = 30.96115307993945

$ws.Cells.Item(20, 2).Value = 13.00119451493854
$ws.Cells.Item(20, 3).Value = 5.180355961625811
$ws.Cells.Item(20, 4).Value = 15.33242102094147
$ws.Cells.Item(20, 5).Value = 16.76369782154143
$ws.Cells.Item(20, 7).Value = 3.692304812720856
$ws.Cells.Item(20, 10).Value = 9.429182338131163
$ws.Cells.Item(20, 11).Value = 12.21263487637867
$ws.Cells.Item(20, 15).Value = 30.95321352553674

$ws.Cells.Item(21, 2).Value = 13.29604260084683
$ws.Cells.Item(21, 3).Value = 5.380969777787831
$ws.Cells.Item(21, 4).Value = 15.44029119315985
$ws.Cells.Item(21, 5).Value = 16.87384561676916
$ws.Cells.Item(21, 7).Value = 3.689555667436582
$ws.Cells.Item(21, 10).Value = 9.428773891073886
$ws.Cells.Item(21, 11).Value = 12.41134672004634
$ws.Cells.Item(21, 15).Value = 30.93384220484709

$ws.Cells.Item(22, 2).Value = 13.48746789434061
$ws.Cells.Item(22, 3).Value = 5.508685765785364
$ws.Cells.Item(22, 4).Value = 15.51330869103314
$ws.Cells.Item(22, 5).Value = 16.94866630475696
$ws.Cells.Item(22, 7).Value = 3.687825171428814
$ws.Cells.Item(22, 10).Value = 9.429754743256542
$ws.Cells.Item(22, 11).Value = 12.54148729401529
$ws.Cells.Item(22, 15).Value = 30.92663364419945

$ws.Cells.Item(23, 2).Value = 13.38543869823746
$ws.Cells.Item(23, 3).Value = 5.440844051452929
$ws.Cells.Item(23, 4).Value = 15.47411283719951
$ws.Cells.Item(23, 5).Value = 16.908479223269
$ws.Cells.Item(23, 7).Value = 3.688742729705655
$ws.Cells.Item(23, 10).Value = 9.429117265055103
$ws.Cells.Item(23, 11).Value = 12.47201723453865
$ws.Cells.Item(23, 15).Value = 30.92998174333018

$ws.Cells.Item(24, 2).Value = 12.99635872942937
$ws.Cells.Item(24, 3).Value = 5.177022487000587
$ws.Cells.Item(24, 4).Value = 15.33070162557256
$ws.Cells.Item(24, 5).Value = 16.7619464993918
$ws.Cells.Item(24, 7).Value = 3.692350868638458
$ws.Cells.Item(24, 10).Value = 9.429210120103123
$ws.Cells.Item(24, 11).Value = 12.2093945777312
$ws.Cells.Item(24, 15).Value = 30.9536220703837

$ws.Cells.Item(25, 2).Value = 12.57345078116374
$ws.Cells.Item(25, 3).Value = 4.879208067861217
$ws.Cells.Item(25, 4).Value = 15.18733895854464
$ws.Cells.Item(25, 5).Value = 16.61656096549197
$ws.Cells.Item(25, 7).Value = 3.696529216973804
$ws.Cells.Item(25, 10).Value = 9.43468443686942
$ws.Cells.Item(25, 11).Value = 11.92862840697348
$ws.Cells.Item(25, 15).Value = 31.00247317119226
